$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2883.3333
$ws.Range("I70").Value = 2466.6667
$ws.Range("J70").Value = 3300
$ws.Range("K70").Value = 7400.000100000001
$ws.Range("L70").Value = 9900
$ws.Range("M70").Value = -7130.000100000001
$ws.Range("N70").Value = -10440
$ws.Range("H73").Value = 2883.3333
$ws.Range("I73").Value = 2466.6667
$ws.Range("J73").Value = 3300
$ws.Range("K73").Value = 7400.000100000001
$ws.Range("L73").Value = 9900
$ws.Range("M73").Value = -6464.000100000001
$ws.Range("N73").Value = -11772
$ws.Range("H135").Value = 8441
$ws.Range("I135").Value = 311
$ws.Range("K135").Value = 2799
$ws.Range("M135").Value = -264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1115304
$ws.Range("I32").Value = 1183945.8
$ws.Range("K32").Value = 1183945.8
$ws.Range("M32").Value = -1183658.8
$ws.Range("H63").Value = 1900
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1900
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102774
$ws.Range("H110").Value = 2070.2222
$ws.Range("I110").Value = 2079
$ws.Range("K110").Value = 2079
$ws.Range("M110").Value = -34
$ws.Range("H112").Value = 27738.666
$ws.Range("J112").Value = 27738.666
$ws.Range("L112").Value = 27738.666
$ws.Range("N112").Value = -30692.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38479.715
$ws.Range("I20").Value = 61984.883
$ws.Range("K20").Value = 61984.883
$ws.Range("M20").Value = -61737.883
$ws.Range("H100").Value = 41999.5
$ws.Range("J100").Value = 41999.5
$ws.Range("L100").Value = 41999.5
$ws.Range("N100").Value = -44163.5
$ws.Range("H105").Value = 1969.9474
$ws.Range("I105").Value = 1674.6072
$ws.Range("K105").Value = 1674.6072
$ws.Range("M105").Value = 72.39280000000008
$ws.Range("H107").Value = 1099.8889
$ws.Range("I107").Value = 1149.8572
$ws.Range("K107").Value = 1149.8572
$ws.Range("M107").Value = 770.1428000000001
$ws.Range("H134").Value = 4635658.5
$ws.Range("I134").Value = 3184.4814
$ws.Range("J134").Value = 18533080
$ws.Range("K134").Value = 9553.4442
$ws.Range("L134").Value = 55599240
$ws.Range("M134").Value = -7018.4442
$ws.Range("N134").Value = -55604310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 62108.633
$ws.Range("I16").Value = 13667.8
$ws.Range("K16").Value = 13667.8
$ws.Range("M16").Value = -13380.8
$ws.Range("H22").Value = 2657.5
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 2762.7273
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2762.7273
$ws.Range("M22").Value = -1150
$ws.Range("N22").Value = -3462.7273
$ws.Range("H57").Value = 24000
$ws.Range("J57").Value = 24000
$ws.Range("L57").Value = 24000
$ws.Range("N57").Value = -25120
$ws.Range("H86").Value = 36875.195
$ws.Range("I86").Value = 59650.94
$ws.Range("K86").Value = 59650.94
$ws.Range("M86").Value = -58527.94
$ws.Range("H89").Value = 36875.195
$ws.Range("I89").Value = 59650.94
$ws.Range("K89").Value = 298254.7
$ws.Range("M89").Value = -292638.7
$ws.Range("H113").Value = 62108.633
$ws.Range("I113").Value = 13667.8
$ws.Range("K113").Value = 13667.8
$ws.Range("M113").Value = -11497.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 684.7857
$ws.Range("J12").Value = 370.66666
$ws.Range("L12").Value = 1111.99998
$ws.Range("N12").Value = -1457.99998
$ws.Range("H37").Value = 69736.84
$ws.Range("J37").Value = 69736.84
$ws.Range("L37").Value = 209210.52
$ws.Range("N37").Value = -209434.52
$ws.Range("H113").Value = 29450.637
$ws.Range("J113").Value = 33995.367
$ws.Range("L113").Value = 101986.101
$ws.Range("N113").Value = -106326.101
$ws.Range("H117").Value = 913.5
$ws.Range("J117").Value = 497.6
$ws.Range("L117").Value = 1492.8
$ws.Range("N117").Value = -8376.799999999999
$ws.Range("I118").Value = 997.5
$ws.Range("J118").Value = 800
$ws.Range("K118").Value = 2992.5
$ws.Range("L118").Value = 2400
$ws.Range("M118").Value = -1749.5
$ws.Range("N118").Value = -4886
$ws.Range("H131").Value = 2849.25
$ws.Range("J131").Value = 3318.25
$ws.Range("L131").Value = 9954.75
$ws.Range("N131").Value = -20034.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2988.5293
$ws.Range("I80").Value = 2797.125
$ws.Range("J80").Value = 3158.6667
$ws.Range("K80").Value = 2797.125
$ws.Range("L80").Value = 3158.6667
$ws.Range("M80").Value = -1799.125
$ws.Range("N80").Value = -5154.6667
$ws.Range("H83").Value = 2988.5293
$ws.Range("I83").Value = 2797.125
$ws.Range("J83").Value = 3158.6667
$ws.Range("K83").Value = 13985.625
$ws.Range("L83").Value = 15793.3335
$ws.Range("M83").Value = -8993.625
$ws.Range("N83").Value = -25777.3335
$ws.Range("H109").Value = 17642
$ws.Range("J109").Value = 17642
$ws.Range("L109").Value = 17642
$ws.Range("N109").Value = -19722
$ws.Range("H122").Value = 6386.1113
$ws.Range("I122").Value = 8743.75
$ws.Range("K122").Value = 26231.25
$ws.Range("M122").Value = -23781.25
$ws.Range("H132").Value = 14608.333
$ws.Range("I132").Value = 16057.625
$ws.Range("J132").Value = 3014
$ws.Range("K132").Value = 48172.875
$ws.Range("L132").Value = 9042
$ws.Range("M132").Value = -45642.875
$ws.Range("N132").Value = -14102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 28867
$ws.Range("J103").Value = 28867
$ws.Range("L103").Value = 28867
$ws.Range("N103").Value = -31211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6946515.5
$ws.Range("I132").Value = 8773657
$ws.Range("K132").Value = 26320971
$ws.Range("M132").Value = -26318441
$ws.Range("H136").Value = 4819697.5
$ws.Range("I136").Value = 2719523
$ws.Range("J136").Value = 10420163
$ws.Range("K136").Value = 8158569
$ws.Range("L136").Value = 31260489
$ws.Range("M136").Value = -8156019
$ws.Range("N136").Value = -31265589
